$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0af57c933cc8c40123709b7fb16e9f578d22c901/e2e/40efff9e-48de-4d28-8aea-d04287e38e3d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b0cb2916c3d46b18eea790756707604c1885339/e2e/40efff9e-48de-4d28-8aea-d04287e38e3d.md."

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Columns.Item(16).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

$ws.Range("I8").Value = "40efff9e-48de-4d28-8aea-d04287e38e3d.md"
$ws.Hyperlinks.Add($ws.Range("I8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b0cb2916c3d46b18eea790756707604c1885339/e2e/40efff9e-48de-4d28-8aea-d04287e38e3d.md", "", "", "40efff9e-48de-4d28-8aea-d04287e38e3d.md")
$ws.Range("J8").Value = "40efff9e-48de-4d28-8aea-d04287e38e3d.8fa37abd8ed77ffba20682d68308b37b595922b4.zh-cn.xlf"
$ws.Range("K8").Value = "2016-08-23 20:45:31"
$ws.Range("P8").Value = $errorDetail

# ---- de-de sheet ----
$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Columns.Item(16).ColumnWidth = $ws2.Columns.Item(1).ColumnWidth

$ws2.Range("I8").Value = "40efff9e-48de-4d28-8aea-d04287e38e3d.md"
$ws2.Hyperlinks.Add($ws2.Range("I8"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b0cb2916c3d46b18eea790756707604c1885339/e2e/40efff9e-48de-4d28-8aea-d04287e38e3d.md", "", "", "40efff9e-48de-4d28-8aea-d04287e38e3d.md")
$ws2.Range("J8").Value = "40efff9e-48de-4d28-8aea-d04287e38e3d.8fa37abd8ed77ffba20682d68308b37b595922b4.de-de.xlf"
$ws2.Range("K8").Value = "2016-08-23 20:45:38"
$ws2.Range("P8").Value = $errorDetail
